# Scheduled runner update: refresh computed market/profit figures on the
# per-job "Profits" sheets (columns H, I, J, K, L, M, N hold average price
# and profit calculations that are recalculated from live market data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2073.1538
$ws.Cells.Item(40, 9).Value = 1979
$ws.Cells.Item(40, 11).Value = 1979
$ws.Cells.Item(40, 13).Value = -1804

$ws.Cells.Item(74, 8).Value = 4000
$ws.Cells.Item(74, 9).Value = 3000
$ws.Cells.Item(74, 11).Value = 3000
$ws.Cells.Item(74, 13).Value = -2064

$ws.Cells.Item(77, 8).Value = 4000
$ws.Cells.Item(77, 9).Value = 3000
$ws.Cells.Item(77, 11).Value = 15000
$ws.Cells.Item(77, 13).Value = -10320

$ws.Cells.Item(93, 8).Value = 88360.4
$ws.Cells.Item(93, 10).Value = 88360.4
$ws.Cells.Item(93, 12).Value = 88360.4
$ws.Cells.Item(93, 14).Value = -93352.4

$ws.Cells.Item(106, 8).Value = 2768.5386
$ws.Cells.Item(106, 9).Value = 3456.889
$ws.Cells.Item(106, 11).Value = 3456.889
$ws.Cells.Item(106, 13).Value = -2825.889

$ws.Cells.Item(107, 8).Value = 855.2778
$ws.Cells.Item(107, 9).Value = 560.7692
$ws.Cells.Item(107, 11).Value = 560.7692
$ws.Cells.Item(107, 13).Value = 1359.2308

$ws.Cells.Item(113, 8).Value = 18420.834
$ws.Cells.Item(113, 9).Value = 31410.8
$ws.Cells.Item(113, 10).Value = 2183.375
$ws.Cells.Item(113, 11).Value = 31410.8
$ws.Cells.Item(113, 12).Value = 2183.375
$ws.Cells.Item(113, 13).Value = -28156.8
$ws.Cells.Item(113, 14).Value = -8691.375

$ws.Cells.Item(132, 8).Value = 809.7292
$ws.Cells.Item(132, 9).Value = 751.11365
$ws.Cells.Item(132, 11).Value = 2253.34095
$ws.Cells.Item(132, 13).Value = 276.6590500000002

$ws.Cells.Item(141, 8).Value = 3501896
$ws.Cells.Item(141, 9).Value = 4000738.2
$ws.Cells.Item(141, 11).Value = 12002214.6
$ws.Cells.Item(141, 13).Value = -11997034.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 758.27026
$ws.Cells.Item(74, 10).Value = 2517.5
$ws.Cells.Item(74, 12).Value = 2517.5
$ws.Cells.Item(74, 14).Value = -4265.5

$ws.Cells.Item(77, 8).Value = 758.27026
$ws.Cells.Item(77, 10).Value = 2517.5
$ws.Cells.Item(77, 12).Value = 12587.5
$ws.Cells.Item(77, 14).Value = -21323.5

$ws.Cells.Item(88, 8).Value = 2532.8667
$ws.Cells.Item(88, 9).Value = 1889
$ws.Cells.Item(88, 10).Value = 3268.7144
$ws.Cells.Item(88, 11).Value = 1889
$ws.Cells.Item(88, 12).Value = 3268.7144
$ws.Cells.Item(88, 13).Value = -1483
$ws.Cells.Item(88, 14).Value = -4080.7144

$ws.Cells.Item(91, 8).Value = 2532.8667
$ws.Cells.Item(91, 9).Value = 1889
$ws.Cells.Item(91, 10).Value = 3268.7144
$ws.Cells.Item(91, 11).Value = 1889
$ws.Cells.Item(91, 12).Value = 3268.7144
$ws.Cells.Item(91, 13).Value = -485
$ws.Cells.Item(91, 14).Value = -6076.7144

$ws.Cells.Item(102, 8).Value = 1321.091
$ws.Cells.Item(102, 9).Value = 1179.4
$ws.Cells.Item(102, 11).Value = 1179.4
$ws.Cells.Item(102, 13).Value = 442.5999999999999

$ws.Cells.Item(130, 8).Value = 44126.375
$ws.Cells.Item(130, 10).Value = 44126.375
$ws.Cells.Item(130, 12).Value = 44126.375
$ws.Cells.Item(130, 14).Value = -54166.375

$ws.Cells.Item(132, 8).Value = 1446.0571
$ws.Cells.Item(132, 9).Value = 1093.5
$ws.Cells.Item(132, 10).Value = 2464.5557
$ws.Cells.Item(132, 11).Value = 3280.5
$ws.Cells.Item(132, 12).Value = 7393.6671
$ws.Cells.Item(132, 13).Value = -750.5
$ws.Cells.Item(132, 14).Value = -12453.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 186412.9
$ws.Cells.Item(86, 9).Value = 5442.75
$ws.Cells.Item(86, 10).Value = 669000
$ws.Cells.Item(86, 11).Value = 5442.75
$ws.Cells.Item(86, 12).Value = 669000
$ws.Cells.Item(86, 13).Value = -4319.75
$ws.Cells.Item(86, 14).Value = -671246

$ws.Cells.Item(89, 8).Value = 186412.9
$ws.Cells.Item(89, 9).Value = 5442.75
$ws.Cells.Item(89, 10).Value = 669000
$ws.Cells.Item(89, 11).Value = 27213.75
$ws.Cells.Item(89, 12).Value = 3345000
$ws.Cells.Item(89, 13).Value = -21597.75
$ws.Cells.Item(89, 14).Value = -3356232

$ws.Cells.Item(99, 8).Value = 1608.1428
$ws.Cells.Item(99, 9).Value = 1452.6
$ws.Cells.Item(99, 10).Value = 1997
$ws.Cells.Item(99, 11).Value = 1452.6
$ws.Cells.Item(99, 12).Value = 1997
$ws.Cells.Item(99, 13).Value = 45.40000000000009
$ws.Cells.Item(99, 14).Value = -4993

$ws.Cells.Item(134, 8).Value = 6800.5454
$ws.Cells.Item(134, 9).Value = 7038
$ws.Cells.Item(134, 10).Value = 1814
$ws.Cells.Item(134, 11).Value = 21114
$ws.Cells.Item(134, 12).Value = 5442
$ws.Cells.Item(134, 13).Value = -18579
$ws.Cells.Item(134, 14).Value = -10512

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 621.7778
$ws.Cells.Item(16, 9).Value = 536.1539
$ws.Cells.Item(16, 11).Value = 536.1539
$ws.Cells.Item(16, 13).Value = -249.1539

$ws.Cells.Item(92, 8).Value = 27999.2
$ws.Cells.Item(92, 10).Value = 27999.2
$ws.Cells.Item(92, 12).Value = 27999.2
$ws.Cells.Item(92, 14).Value = -32991.2

$ws.Cells.Item(95, 8).Value = 24249.75
$ws.Cells.Item(95, 10).Value = 24249.75
$ws.Cells.Item(95, 12).Value = 24249.75
$ws.Cells.Item(95, 14).Value = -29741.75

$ws.Cells.Item(113, 8).Value = 621.7778
$ws.Cells.Item(113, 9).Value = 536.1539
$ws.Cells.Item(113, 11).Value = 536.1539
$ws.Cells.Item(113, 13).Value = 1633.8461

$ws.Cells.Item(122, 8).Value = 4144.7144
$ws.Cells.Item(122, 9).Value = 2250
$ws.Cells.Item(122, 11).Value = 6750
$ws.Cells.Item(122, 13).Value = -4300

$ws.Cells.Item(134, 8).Value = 1914.7646
$ws.Cells.Item(134, 9).Value = 1706.5807
$ws.Cells.Item(134, 11).Value = 5119.742099999999
$ws.Cells.Item(134, 13).Value = -2584.742099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 834.3333
$ws.Cells.Item(122, 9).Value = 635
$ws.Cells.Item(122, 11).Value = 5715
$ws.Cells.Item(122, 13).Value = -3265

$ws.Cells.Item(123, 8).Value = 833.3333
$ws.Cells.Item(123, 9).Value = 833.3333
$ws.Cells.Item(123, 11).Value = 2499.9999
$ws.Cells.Item(123, 13).Value = -49.9998999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5122
$ws.Cells.Item(70, 9).Value = 4800
$ws.Cells.Item(70, 11).Value = 4800
$ws.Cells.Item(70, 13).Value = -4530

$ws.Cells.Item(73, 8).Value = 5122
$ws.Cells.Item(73, 9).Value = 4800
$ws.Cells.Item(73, 11).Value = 4800
$ws.Cells.Item(73, 13).Value = -3864

$ws.Cells.Item(113, 8).Value = 747.5
$ws.Cells.Item(113, 9).Value = 543.4286
$ws.Cells.Item(113, 11).Value = 543.4286
$ws.Cells.Item(113, 13).Value = 1626.5714

$ws.Cells.Item(132, 8).Value = 1481540.2
$ws.Cells.Item(132, 9).Value = 1924612.4
$ws.Cells.Item(132, 11).Value = 5773837.199999999
$ws.Cells.Item(132, 13).Value = -5771307.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1233.3334
$ws.Cells.Item(100, 9).Value = 1300
$ws.Cells.Item(100, 11).Value = 1300
$ws.Cells.Item(100, 13).Value = -759

$ws.Cells.Item(136, 8).Value = 1570.0278
$ws.Cells.Item(136, 9).Value = 1228.6552
$ws.Cells.Item(136, 11).Value = 3685.9656
$ws.Cells.Item(136, 13).Value = -1135.9656

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1390.8572
$ws.Cells.Item(100, 9).Value = 1117
$ws.Cells.Item(100, 10).Value = 1756
$ws.Cells.Item(100, 11).Value = 2234
$ws.Cells.Item(100, 12).Value = 3512
$ws.Cells.Item(100, 13).Value = -1693
$ws.Cells.Item(100, 14).Value = -4594

$ws.Cells.Item(107, 8).Value = 646.85187
$ws.Cells.Item(107, 9).Value = 476.94446
$ws.Cells.Item(107, 10).Value = 986.6667
$ws.Cells.Item(107, 11).Value = 1430.83338
$ws.Cells.Item(107, 12).Value = 2960.0001
$ws.Cells.Item(107, 13).Value = 489.16662
$ws.Cells.Item(107, 14).Value = -6800.0001

$ws.Cells.Item(136, 8).Value = 30866102
$ws.Cells.Item(136, 9).Value = 55557130
$ws.Cells.Item(136, 11).Value = 166671390
$ws.Cells.Item(136, 13).Value = -166668840
